$d = $word.ActiveDocument
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$t = $d.Tables.Item(15)
$t.Delete()
$paraB = $d.Range(2609, 2610)
$paraB.Delete()

# Try including the paragraph mark itself (extend range by 1 more char)
$paraCFull = $d.Range(2609, 2619)
Write-Host "paraCFull text: [" $paraCFull.Text "]"
$paraCFull.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:bookmarkStart w:id='77' w:name='_GoBack'/><w:bookmarkEnd w:id='77'/></w:p>")
Write-Host "Content end:" $d.Content.End
